$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update two existing values that were re-computed (D22 and C23)
$ws.Range("D22").Value = 0.7115302104241067
$ws.Range("C23").Value = 0.3054124294241067

# Add a new data row (row 24) for the newest ifoCAST sampling date,
# matching the label style used by the other date rows (A2:A23)
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A24").Value = "2025-09-04_diff"
$ws.Range("B24").Value = 0.05603945542410671

Write-Output "done"
